# Refresh the cryptos price/volume table (GitHub Actions daily snapshot update).
# Note: several "Price" (column D) values are plain decimal-looking strings
# (e.g. "249.07"). Assigning those straight to Range.Value lets Excel's COM
# layer auto-coerce them to numeric cells, which would drop the original
# text formatting used throughout column D (e.g. "37.118.83", "2.049.20").
# To keep them as text - matching the rest of the sheet - we briefly force
# the cell to a text number format before the assignment, then restore the
# default "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.118.83'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '2.049.20'
$ws.Range("E3").Value = '  -0.57%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.79'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +8.02%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0793'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.08'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.62%  '
$ws.Range("D13").Value = '2.349.22'
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.832'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +10.14%  '
$ws.Range("D16").Value = '2.050.92'
$ws.Range("E16").Value = '  -0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '18.31'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +28.33%  '
$ws.Range("D18").Value = '37.094.06'
$ws.Range("E18").Value = '  +0.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '76.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.89%  '
$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -4.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("E25").Value = '  +10.29%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.44'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.51%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.23'
$ws.Range("D28").Style = "Normal"
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.84'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.25%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0630'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.62'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0884'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  -0.09%  '
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.75'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.51%  '
$ws.Range("E38").Value = '  +3.09%  '
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("E40").Value = '  +11.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.16'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +20.59%  '
$ws.Range("E42").Value = '  -0.40%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.33%  '
$ws.Range("E44").Value = '  -0.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '97.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.94%  '
$ws.Range("D48").Value = '1.293.53'
$ws.Range("E48").Value = '  -0.67%  '
$ws.Range("E49").Value = '  -1.23%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.31%  '
$ws.Range("D51").Value = '2.237.98'
$ws.Range("E51").Value = '  -0.55%  '
